$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.018300890922546
$ws.Range("B1").Value = 1.748725891113281
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.004086971282959
$ws.Range("E1").Value = 1.261240482330322
